# Correction on heat sector data (1)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 1: add a "PJ" unit header merged across B1:D1 (A1 already says "PJ") ---
$ws.Cells.Item(1, 2).Font.Name = "Arial"
$ws.Cells.Item(1, 2).Font.Size = 9
$ws.Cells.Item(1, 2).Font.Color = 6710886
$ws.Cells.Item(1, 2).HorizontalAlignment = -4108
$ws.Range("B1:D1").Merge()

# --- Row 2: rename the "Electric boilers" column header to "Null" ---
$ws.Range("B2").Value = "Null"

# --- Rows 3-33: corrected data values (B=Electric boilers/Null col, C=Heat pump - air, D=Heat pump - ground) ---
$ws.Range("B3").Value = 25.881
$ws.Range("C3").Value = 29.185
$ws.Range("D3").Value = 58.816

$ws.Range("B4").Value = 40.74
$ws.Range("C4").Value = 81.743
$ws.Range("D4").Value = 55.444

$ws.Range("B5").Value = 19.926
$ws.Range("C5").Value = 4.217

$ws.Range("B6").Value = 39.348
$ws.Range("C6").Value = 74.328
$ws.Range("D6").Value = 51.34

$ws.Range("B7").Value = 4.052
$ws.Range("C7").Value = 3.781
$ws.Range("D7").Value = 0.005

$ws.Range("B8").Value = 41.154
$ws.Range("C8").Value = 52.127
$ws.Range("D8").Value = 53.379

$ws.Range("B9").Value = 147.8
$ws.Range("C9").Value = 601.904
$ws.Range("D9").Value = 402.276

$ws.Range("B10").Value = 18.726
$ws.Range("C10").Value = 14.544
$ws.Range("D10").Value = 36.948

$ws.Range("B11").Value = 6.007
$ws.Range("C11").Value = 4.118
$ws.Range("D11").Value = 5.535

$ws.Range("B12").Value = 39.061
$ws.Range("C12").Value = 22.74
$ws.Range("D12").Value = 29.965

$ws.Range("B13").Value = 157.574
$ws.Range("C13").Value = 84.171
$ws.Range("D13").Value = 29.714

$ws.Range("B14").Value = 55.402
$ws.Range("C14").Value = 31.181
$ws.Range("D14").Value = 54.203

$ws.Range("B15").Value = 356.988
$ws.Range("C15").Value = 311.838
$ws.Range("D15").Value = 277.94

$ws.Range("B16").Value = 10.431
$ws.Range("C16").Value = 7.779
$ws.Range("D16").Value = 18.125

$ws.Range("B17").Value = 20.085
$ws.Range("C17").Value = 45.664
$ws.Range("D17").Value = 38.149

$ws.Range("B18").Value = 16.806
$ws.Range("C18").Value = 28.421
$ws.Range("D18").Value = 15.12

$ws.Range("B19").Value = 2.338
$ws.Range("C19").Value = 3.509
$ws.Range("D19").Value = 2.716

$ws.Range("B20").Value = 191.726
$ws.Range("C20").Value = 358.313
$ws.Range("D20").Value = 171.271

$ws.Range("B21").Value = 5.301
$ws.Range("C21").Value = 3.425
$ws.Range("D21").Value = 4.311

$ws.Range("B22").Value = 3.371
$ws.Range("C22").Value = 7.804
$ws.Range("D22").Value = 5.222

$ws.Range("B23").Value = 4.115
$ws.Range("C23").Value = 5.365
$ws.Range("D23").Value = 5.892

$ws.Range("B24").Value = 1.078
$ws.Range("C24").Value = 0.004

$ws.Range("B25").Value = 80.632
$ws.Range("C25").Value = 227.669
$ws.Range("D25").Value = 53.266

$ws.Range("B26").Value = 94.754
$ws.Range("C26").Value = 32.673
$ws.Range("D26").Value = 2.501

$ws.Range("B27").Value = 88.579
$ws.Range("C27").Value = 85.207
$ws.Range("D27").Value = 120.497

$ws.Range("B28").Value = 23.205
$ws.Range("C28").Value = 6.202
$ws.Range("D28").Value = 0.013

$ws.Range("B29").Value = 39.132
$ws.Range("C29").Value = 16.979
$ws.Range("D29").Value = 45.208

$ws.Range("B30").Value = 127.128
$ws.Range("C30").Value = 10.817
$ws.Range("D30").Value = 34.825

$ws.Range("B31").Value = 7.802
$ws.Range("C31").Value = 5.716
$ws.Range("D31").Value = 7.356

$ws.Range("B32").Value = 14.42
$ws.Range("C32").Value = 24.129
$ws.Range("D32").Value = 15.145

$ws.Range("B33").Value = 401.32
$ws.Range("C33").Value = 364.984
$ws.Range("D33").Value = 206.912

# --- Column widths (best fit on the now-narrower Region / label columns) ---
$ws.Columns("A").ColumnWidth = 6
$ws.Columns("B").ColumnWidth = 5

# --- Selection cosmetics ---
$ws.Range("F7").Select()
